# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Feria Lagunitas de Puerto Montt - Mango"
# as row 235, pushing every following row down by one (old A1:T297 -> A1:T298).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 235:297 down to 236:298, leaving a blank row 235 to fill in.
$ws.Rows.Item(235).Insert()

# Populate the newly inserted row 235 with the new data point.
$ws.Cells.Item(235, 1).Value = 4
$ws.Cells.Item(235, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(235, 3).Value = "Los Lagos"
$ws.Cells.Item(235, 4).Value = 44932
$ws.Cells.Item(235, 5).Value = 10
$ws.Cells.Item(235, 6).Value = "Fruta"
$ws.Cells.Item(235, 7).Value = 100108
$ws.Cells.Item(235, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(235, 9).Value = 100108002
$ws.Cells.Item(235, 10).Value = "Mango"
$ws.Cells.Item(235, 11).Value = "Sin especificar"
$ws.Cells.Item(235, 12).Value = "Primera"
$ws.Cells.Item(235, 13).Value = 160
$ws.Cells.Item(235, 14).Value = 8000
$ws.Cells.Item(235, 15).Value = 8500
$ws.Cells.Item(235, 16).Value = 8250
$ws.Cells.Item(235, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(235, 18).Value = "Brasil"
$ws.Cells.Item(235, 19).Value = 2062
$ws.Cells.Item(235, 20).Value = 4
